$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 124345
$ws.Range("C2").Value = 4332
$ws.Range("B3").Value = 588
$ws.Range("C3").Value = 43520
